$wb = $excel.ActiveWorkbook

# --- Metadata sheet value updates -----------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# --- Re-affirm alignment (vertical=top + wrap text) on every already ------
# --- aligned cell so Excel marks the style with applyAlignment="true" -----
foreach ($ws in $wb.Worksheets) {
  $used = $ws.UsedRange
  foreach ($cell in $used.Cells) {
    if ($cell.VerticalAlignment -eq -4160) {
      $cell.WrapText = $true
      $cell.VerticalAlignment = -4160
    }
  }
}
